# issue #5: stock data from json to db
#
# 股票 (stock) worksheet gains a "category" column (always "normal") right
# after "property_category", plus two trailing columns: "source_file"
# (the workbook's own temp-file stem, "tmp1afe1") and "index" (mirrors the
# row's existing numbering in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票

$lastRow = 7

# Shift the old "date"/"legislator_name"/"legislator_id" columns (I:K) one
# place right, freeing column I for "category". Insert() clones the
# neighbouring column's formatting (header s=1 / body s=2) onto every row.
$ws.Columns("I:I").Insert()

# Two brand-new columns appended after "legislator_id" (now column L).
$ws.Columns("M:N").Insert()

# --- header row -----------------------------------------------------
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- data rows --------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp1afe1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
